$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new columns before column D (shifts D:K -> F:M)
$ws.Range("D1:E1").EntireColumn.Insert()

# Copy formatting from column F (old D, now shifted) into new D:E columns
$ws.Columns("F").Copy()
$ws.Columns("D:E").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Set new quarter-end dates in row 7/38/80 (2018-12-31 and 2018-09-30)
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373

# Set new financial data values for columns D and E across all data rows
$ws.Range("D8").Value = 17400
$ws.Range("E8").Value = 6900
$ws.Range("D9").Value = 600
$ws.Range("E9").Value = 600
$ws.Range("D10").Value = 16800
$ws.Range("E10").Value = 6300
$ws.Range("D12").Value = 12100
$ws.Range("E12").Value = 13800
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = "NA"
$ws.Range("E14").Value = "NA"
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 29200
$ws.Range("E17").Value = 29900
$ws.Range("D18").Value = -11800
$ws.Range("E18").Value = -23000
$ws.Range("D20").Value = 800
$ws.Range("E20").Value = 700
$ws.Range("D21").Value = -10100
$ws.Range("E21").Value = -21300
$ws.Range("D22").Value = 5200
$ws.Range("E22").Value = 5300
$ws.Range("D23").Value = -16300
$ws.Range("E23").Value = -27500
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = -16300
$ws.Range("E26").Value = -27500
$ws.Range("D27").Value = -16300
$ws.Range("E27").Value = -27500
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -800
$ws.Range("E32").Value = -700
$ws.Range("D33").Value = -16300
$ws.Range("E33").Value = -27500
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = -16300
$ws.Range("E35").Value = -27500
$ws.Range("D41").Value = 80400
$ws.Range("E41").Value = 48900
$ws.Range("D42").Value = 79700
$ws.Range("E42").Value = 138400
$ws.Range("D43").Value = 5900
$ws.Range("E43").Value = 2900
$ws.Range("D44").Value = 4700
$ws.Range("E44").Value = 4800
$ws.Range("D45").Value = 2700
$ws.Range("E45").Value = 3400
$ws.Range("D46").Value = 173300
$ws.Range("E46").Value = 198400
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 15900
$ws.Range("E48").Value = 16300
$ws.Range("D49").Value = 94700
$ws.Range("E49").Value = 95100
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 300
$ws.Range("E52").Value = 400
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 284100
$ws.Range("E54").Value = 310200
$ws.Range("D57").Value = 17800
$ws.Range("E57").Value = 51500
$ws.Range("D58").Value = 1100
$ws.Range("E58").Value = 1100
$ws.Range("D59").Value = 17900
$ws.Range("E59").Value = 15900
$ws.Range("D60").Value = 36800
$ws.Range("E60").Value = 68500
$ws.Range("D61").Value = 243900
$ws.Range("E61").Value = 243800
$ws.Range("D62").Value = 29900
$ws.Range("E62").Value = 27300
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 310500
$ws.Range("E66").Value = 339600
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -1471600
$ws.Range("E72").Value = -1471500
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = -26400
$ws.Range("E76").Value = -29400
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D81").Value = -16300
$ws.Range("E81").Value = -27500
$ws.Range("D83").Value = 900
$ws.Range("E83").Value = 900
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = -27100
$ws.Range("E89").Value = -22200
$ws.Range("D91").Value = 0
$ws.Range("E91").Value = 0
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = 58800
$ws.Range("E94").Value = 56400
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -100
$ws.Range("E100").Value = -400
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = 31500
$ws.Range("E102").Value = 33900

# Apply corrections to restated figures for the quarter ending 2017-12-31 (column H after insert)
$ws.Range("H8").Value = 34400
$ws.Range("H10").Value = 33900
$ws.Range("H12").Value = 42600
$ws.Range("H17").Value = 59100
$ws.Range("H18").Value = -24700
$ws.Range("H21").Value = -23200
$ws.Range("H23").Value = -26300
$ws.Range("H26").Value = -26300
$ws.Range("H27").Value = -26300
$ws.Range("H33").Value = -22300
$ws.Range("H35").Value = -22300
$ws.Range("H48").Value = 94000
$ws.Range("H57").Value = 38800
$ws.Range("H59").Value = 52600
$ws.Range("H60").Value = 105500
$ws.Range("H62").Value = 31200
$ws.Range("H66").Value = 368300
$ws.Range("H72").Value = -1365200
$ws.Range("H76").Value = 68300
$ws.Range("H81").Value = -22300
